# Update the "when was the transfer first discussed" question text:
# "Date the transfer was first discussed with the incoming or outgoing trust"
# becomes
# "Date the transfer was first discussed with a trust"
$d = $word.ActiveDocument

$find = $d.Content.Find
$found = $find.Execute("the incoming or outgoing trust", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "a trust", 2)

if (-not $found) {
    throw "Could not find the text to replace."
}
